$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.269.10'
$ws.Range('E2').Value = '  -0.63%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.384.33'
$ws.Range('E3').Value = '  -3.90%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '549.17'
$ws.Range('E5').Value = '  -1.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.62'
$ws.Range('E6').Value = '  -4.30%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.531'
$ws.Range('E8').Value = '  -11.52%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.384.07'
$ws.Range('E9').Value = '  -3.86%  '
$ws.Range('E10').Value = '  -2.68%  '
$ws.Range('E11').Value = '  +0.22%  '
$ws.Range('E12').Value = '  -3.34%  '
$ws.Range('E13').Value = '  -3.05%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.39'
$ws.Range('E14').Value = '  -4.16%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.820.54'
$ws.Range('E15').Value = '  -3.80%  '
$ws.Range('E16').Value = '  -2.47%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '60.840.67'
$ws.Range('E17').Value = '  -1.18%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.387.23'
$ws.Range('E18').Value = '  -4.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.74'
$ws.Range('E19').Value = '  -4.66%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.12'
$ws.Range('E20').Value = '  -2.54%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '318.61'
$ws.Range('E21').Value = '  -1.21%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.69'
$ws.Range('E22').Value = '  -6.75%  '
$ws.Range('E23').Value = '  -0.11%  '
$ws.Range('E24').Value = '  +0.47%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '63.79'
$ws.Range('E25').Value = '  -1.01%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.20'
$ws.Range('E26').Value = '  +4.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.507.44'
$ws.Range('E28').Value = '  -3.79%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '529.44'
$ws.Range('E29').Value = '  -6.20%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0926'
$ws.Range('E30').Value = '  -9.02%  '
$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.43'
$ws.Range('E31').Value = '  -5.73%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.09'
$ws.Range('E32').Value = '  -3.25%  '
$ws.Range('E33').Value = '  -3.29%  '
$ws.Range('E34').Value = '  -5.08%  '
$ws.Range('E35').Value = '  -0.83%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.51'
$ws.Range('E37').Value = '  -7.15%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.67'
$ws.Range('E38').Value = '  -5.46%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.375'
$ws.Range('E39').Value = '  -2.42%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.86'
$ws.Range('E40').Value = '  +4.91%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.10'
$ws.Range('E41').Value = '  -2.98%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '140.15'
$ws.Range('E42').Value = '  -4.58%  '
$ws.Range('E43').Value = '  -0.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '40.53'
$ws.Range('E44').Value = '  -0.12%  '
$ws.Range('E45').Value = '  -10.99%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '140.95'
$ws.Range('E46').Value = '  -5.64%  '
$ws.Range('B47').Value = 'Filecoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.61'
$ws.Range('E47').Value = '  -1.53%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '20.11'
$ws.Range('E48').Value = '  -9.33%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0519'
$ws.Range('E49').Value = '  -4.60%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.576'
$ws.Range('E50').Value = '  -3.70%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0907'
$ws.Range('E51').Value = '  -4.19%  '
